$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 74, pushing the existing rows 74-92 down to 75-93.
$ws.Rows.Item(74).Insert()

# Populate the newly inserted row 74 with the new weekly price record.
$ws.Cells.Item(74, 1).Value = 4
$ws.Cells.Item(74, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(74, 3).Value = "Los Lagos"
$ws.Cells.Item(74, 4).Value = 45258
$ws.Cells.Item(74, 5).Value = 10
$ws.Cells.Item(74, 6).Value = 300000000
$ws.Cells.Item(74, 7).Value = "Espárragos"
$ws.Cells.Item(74, 8).Value = "Sin especificar"
$ws.Cells.Item(74, 9).Value = "Primera"
$ws.Cells.Item(74, 10).Value = 500
$ws.Cells.Item(74, 11).Value = 2000
$ws.Cells.Item(74, 12).Value = 2000
$ws.Cells.Item(74, 13).Value = 2000
$ws.Cells.Item(74, 14).Value = "$/kilo"
$ws.Cells.Item(74, 15).Value = "Provincia de Linares"
$ws.Cells.Item(74, 16).Value = 2000
$ws.Cells.Item(74, 17).Value = 1
$ws.Cells.Item(74, 18).Value = "Hortaliza"
